$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item(1)
$zhcn     = $wb.Worksheets.Item(2)
$dede     = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------------
# 1. Status text change: "Ready for handoff" -> "Handed back: in sync with en-US"
#    Appears on the Overview sheet (E2,F2,E3,F3) and on each language sheet's
#    Status column (C2,C3).
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$overview.Range("E2").Value2 = $newStatus
$overview.Range("F2").Value2 = $newStatus
$overview.Range("E3").Value2 = $newStatus
$overview.Range("F3").Value2 = $newStatus

$zhcn.Range("C2").Value2 = $newStatus
$zhcn.Range("C3").Value2 = $newStatus

$dede.Range("C2").Value2 = $newStatus
$dede.Range("C3").Value2 = $newStatus

# ---------------------------------------------------------------------------
# 2. Populate "Latest Target File" (I) / "Latest Handback File" (J) /
#    "Latest Handback DateTime" (K) for both rows on both language sheets,
#    mirroring the handoff info now that the handback report has run.
# ---------------------------------------------------------------------------

$row2Doc = "26763ceb-189b-4453-9c18-82c4240f751d.md"
$row3Doc = "e66f312d-5ff1-4418-a3d5-b0ad03944bf8.md"

$row2Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5f6f98fdcae8333a6862f47a2f01f82c7dddc3cb/e2e/26763ceb-189b-4453-9c18-82c4240f751d.md"
$row3Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5f6f98fdcae8333a6862f47a2f01f82c7dddc3cb/e2e/e66f312d-5ff1-4418-a3d5-b0ad03944bf8.md"

# --- zh-cn sheet ---
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $row2Url, "", "", $row2Doc) | Out-Null
$zhcn.Range("J2").Value2 = "26763ceb-189b-4453-9c18-82c4240f751d.a904af15abd2999f0ae13f8757cb336731056cfd.zh-cn.xlf"
$zhcn.Range("K2").Value2 = "2016-08-29 07:03:02"

$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $row3Url, "", "", $row3Doc) | Out-Null
$zhcn.Range("J3").Value2 = "e66f312d-5ff1-4418-a3d5-b0ad03944bf8.45594e1b5742ab6b709b702600ff7b25aca950c0.zh-cn.xlf"
$zhcn.Range("K3").Value2 = "2016-08-29 07:03:02"

# --- de-de sheet ---
$dede.Hyperlinks.Add($dede.Range("I2"), $row2Url, "", "", $row2Doc) | Out-Null
$dede.Range("J2").Value2 = "26763ceb-189b-4453-9c18-82c4240f751d.a904af15abd2999f0ae13f8757cb336731056cfd.de-de.xlf"
$dede.Range("K2").Value2 = "2016-08-29 07:03:15"

$dede.Hyperlinks.Add($dede.Range("I3"), $row3Url, "", "", $row3Doc) | Out-Null
$dede.Range("J3").Value2 = "e66f312d-5ff1-4418-a3d5-b0ad03944bf8.45594e1b5742ab6b709b702600ff7b25aca950c0.de-de.xlf"
$dede.Range("K3").Value2 = "2016-08-29 07:03:15"

# ---------------------------------------------------------------------------
# 3. Widen columns that now hold longer content (status text, hyperlinked
#    file names, handback xlf names) so everything stays readable.
# ---------------------------------------------------------------------------
$overview.Columns.Item(5).ColumnWidth = 29.9777047293527
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527

$zhcn.Columns.Item(3).ColumnWidth = 29.9777047293527
$zhcn.Columns.Item(9).ColumnWidth = 40
$zhcn.Columns.Item(10).ColumnWidth = 40

$dede.Columns.Item(3).ColumnWidth = 29.9777047293527
$dede.Columns.Item(9).ColumnWidth = 40
$dede.Columns.Item(10).ColumnWidth = 40
